$wb = $excel.ActiveWorkbook

# Identify the existing sheets (before insertion).
$wsBrix  = $wb.Worksheets.Item(2)   # Brix_Daten

# Insert a brand-new worksheet just before "Brix_Daten" -> it becomes the
# 2nd sheet, matching the new sheet order: Startseite, Carport_Liste,
# Brix_Daten, GeländerEdelstahl_Daten.
$wsCarport = $wb.Worksheets.Add($wsBrix)
$wsCarport.Name = "Carport_Liste"

# NOTE: worksheet object references behave like live, index-bound handles
# in this object model, so after inserting a sheet the old $wsBrix handle
# now actually points at the freshly inserted sheet (same index slot).
# Re-fetch every sheet reference we still need by name to make sure they
# point at the right worksheet.
$wsStart   = $wb.Worksheets.Item("Startseite")
$wsBrix    = $wb.Worksheets.Item("Brix_Daten")
$wsCarport = $wb.Worksheets.Item("Carport_Liste")

# --- Populate formatting first (so the later Value assignments keep the
# shared-string / style bookkeeping correct), re-using the same visual
# styles already used on the "Startseite" sheet: row 1 uses the bold
# header style, rows 2-3 use the regular bordered style, row 4 is plain.
$wsStart.Range("A1:B1").Copy()
$wsCarport.Range("A1:C1").PasteSpecial(-4122)   # xlPasteFormats

$wsStart.Range("A2:B2").Copy()
$wsCarport.Range("A2:C3").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

$wsCarport.Rows(1).RowHeight = 30.75
$wsCarport.Rows(2).RowHeight = 15.75
$wsCarport.Rows(3).RowHeight = 15.75

# --- Fill in the table contents.
$wsCarport.Range("A1").Value = "Bezeichnung"
$wsCarport.Range("B1").Value = "Einheit"
$wsCarport.Range("C1").Value = "Preis"

$wsCarport.Range("A2").Value = "Carport "
$wsCarport.Range("B2").Value = "Qm"
$wsCarport.Range("C2").Value = 500

$wsCarport.Range("A3").Value = "Steher "
$wsCarport.Range("B3").Value = "Stk"
$wsCarport.Range("C3").Value = 100

$wsCarport.Range("A4").Value = "Sanwich eindeckung"
$wsCarport.Range("B4").Value = "qm"
$wsCarport.Range("C4").Value = 45

# --- Restore/replicate the selections recorded in each sheet view.
[void]$wsStart.Range("F32").Select()
[void]$wsBrix.Range("A1:C3").Select()

# Carport_Liste is selected last so it ends up as the active tab, matching
# the workbook's activeTab pointing at it.
[void]$wsCarport.Range("C5").Select()
$wsCarport.Activate()
